# Fixed some code error
# Regenerate the "Email Address" column (E) for each student row so that it
# is derived from the Student Name (column C): first letter of the surname
# (lowercase) + last given/middle name token (lowercase) + "@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $nameCell = $ws.Cells.Item($r, 3)
    $name = $nameCell.Value2

    if ([string]::IsNullOrEmpty($name)) {
        continue
    }

    $trimmed = $name.Trim()
    $initial = $trimmed.Substring(0, 1).ToLower()

    $parts = $trimmed -split "\s+"
    $lastWord = $parts[$parts.Length - 1].ToLower()

    $newEmail = $initial + $lastWord + "@gmail.com"

    $ws.Cells.Item($r, 5).Value = $newEmail
}
